# "Generate Report for Handoff"
# The localization status report has a new row result: the entry for
# f51aae06-b332-4504-88e2-84341bc8abb1.md has moved from "In Translation"
# to "Ready for handoff" (a fresh handoff .xlf was generated), so every
# sheet that tracks that file needs its Status / Priority / handoff
# timestamp columns refreshed.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-24 04:13:34"

# ---- zh-cn sheet ---------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-24 04:13:30"

# ---- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-24 04:13:34"
